$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 386 — this shifts rows 386..437 down to 387..438,
# preserving every existing record's data and formatting.
$ws.Rows.Item(386).Insert()

# Populate the newly inserted row 386 with the new weekly record.
$ws.Cells.Item(386, 1).Value = 5
$ws.Cells.Item(386, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(386, 3).Value = "Maule"
$ws.Cells.Item(386, 4).Value = 45077
$ws.Cells.Item(386, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(386, 5).Value = 7
$ws.Cells.Item(386, 6).Value = 100112009
$ws.Cells.Item(386, 7).Value = "Acelga"
$ws.Cells.Item(386, 8).Value = "Sin especificar"
$ws.Cells.Item(386, 9).Value = "Primera"
$ws.Cells.Item(386, 10).Value = 500
$ws.Cells.Item(386, 11).Value = 2000
$ws.Cells.Item(386, 12).Value = 2000
$ws.Cells.Item(386, 13).Value = 2000
$ws.Cells.Item(386, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(386, 15).Value = "Región del Maule"
$ws.Cells.Item(386, 16).Value = 500
$ws.Cells.Item(386, 17).Value = 4
$ws.Cells.Item(386, 18).Value = "Hortaliza"
